# Updates the cryptos price/volume snapshot (GitHub Actions refresh).
# D-column price cells get NumberFormat="@" before the write so numeric-
# looking text (e.g. "0.414") is stored as text rather than coerced to a
# number, then Style is reset to "Normal" so no lasting cell formatting
# change is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '92.044.30'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.89%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.271.68'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '628.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.414'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +9.04%  '
$ws.Range('E8').Value = '  +4.61%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.270.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.25%  '
$ws.Range('E11').Value = '  +1.51%  '
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.874.91'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.803.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.277.93'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.03'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '438.89'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0000195'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +49.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.452.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '77.05'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.76'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '553.87'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.30%  '
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.64'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +23.60%  '
$ws.Range('E37').Value = '  -8.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.64'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.42'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.130'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.997'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.394'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '149.61'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.40%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '180.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.18'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('E48').Value = '  +5.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.28'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('E51').Value = '  +1.80%  '
